$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handoff
# Two new localized files handed off:
#   56ad0c0d-4b68-4968-890c-943c2590523c.md
#   953e0f08-aecb-4b53-a766-72bc9d9918db.md
# Each gets a new row (row 6 / row 7) appended to the three report tables:
#   Overview, zh-cn, de-de
# ---------------------------------------------------------------------------

# ---------------------- Overview sheet --------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$loOverview.ListRows.Add()
$wsOverview.Range("A6").Value = "56ad0c0d-4b68-4968-890c-943c2590523c.md"
$wsOverview.Range("B6").Value = "e2e\56ad0c0d-4b68-4968-890c-943c2590523c.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2dabc17b3f6fd8472fd607b2ff25a6c26e569a10/e2e/56ad0c0d-4b68-4968-890c-943c2590523c.md", "", "", "e2e\56ad0c0d-4b68-4968-890c-943c2590523c.md")
$wsOverview.Range("C6").Value = ".md"
$wsOverview.Range("D6").Value = ""
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-09-05 02:46:36"

$loOverview.ListRows.Add()
$wsOverview.Range("A7").Value = "953e0f08-aecb-4b53-a766-72bc9d9918db.md"
$wsOverview.Range("B7").Value = "e2e\953e0f08-aecb-4b53-a766-72bc9d9918db.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c21cabe021da961c62ae859f803b17101d2a2c0e/e2e/953e0f08-aecb-4b53-a766-72bc9d9918db.md", "", "", "e2e\953e0f08-aecb-4b53-a766-72bc9d9918db.md")
$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("D7").Value = ""
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-09-05 02:46:36"

# ---------------------- zh-cn sheet ------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)

$loZhCn.ListRows.Add()
$wsZhCn.Range("A6").Value = "56ad0c0d-4b68-4968-890c-943c2590523c.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2dabc17b3f6fd8472fd607b2ff25a6c26e569a10/e2e/56ad0c0d-4b68-4968-890c-943c2590523c.md", "", "", "56ad0c0d-4b68-4968-890c-943c2590523c.md")
$wsZhCn.Range("B6").Value = ".md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "e2e"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("F6").Value = "False"
$wsZhCn.Range("G6").Value = "56ad0c0d-4b68-4968-890c-943c2590523c.2dabc17b3f6fd8472fd607b2ff25a6c26e569a10.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "2016-09-05 02:46:31"
$wsZhCn.Range("I6").Value = ""
$wsZhCn.Range("J6").Value = ""
$wsZhCn.Range("K6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L6").Value = ""
$wsZhCn.Range("M6").Value = "True"
$wsZhCn.Range("N6").Value = ""
$wsZhCn.Range("O6").Value = "False"
$wsZhCn.Range("P6").Value = ""

$loZhCn.ListRows.Add()
$wsZhCn.Range("A7").Value = "953e0f08-aecb-4b53-a766-72bc9d9918db.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c21cabe021da961c62ae859f803b17101d2a2c0e/e2e/953e0f08-aecb-4b53-a766-72bc9d9918db.md", "", "", "953e0f08-aecb-4b53-a766-72bc9d9918db.md")
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e2e"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("F7").Value = "False"
$wsZhCn.Range("G7").Value = "953e0f08-aecb-4b53-a766-72bc9d9918db.c21cabe021da961c62ae859f803b17101d2a2c0e.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-09-05 02:46:31"
$wsZhCn.Range("I7").Value = ""
$wsZhCn.Range("J7").Value = ""
$wsZhCn.Range("K7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L7").Value = ""
$wsZhCn.Range("M7").Value = "True"
$wsZhCn.Range("N7").Value = ""
$wsZhCn.Range("O7").Value = "False"
$wsZhCn.Range("P7").Value = ""

# ---------------------- de-de sheet ------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)

$loDeDe.ListRows.Add()
$wsDeDe.Range("A6").Value = "56ad0c0d-4b68-4968-890c-943c2590523c.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2dabc17b3f6fd8472fd607b2ff25a6c26e569a10/e2e/56ad0c0d-4b68-4968-890c-943c2590523c.md", "", "", "56ad0c0d-4b68-4968-890c-943c2590523c.md")
$wsDeDe.Range("B6").Value = ".md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "e2e"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("F6").Value = "False"
$wsDeDe.Range("G6").Value = "56ad0c0d-4b68-4968-890c-943c2590523c.2dabc17b3f6fd8472fd607b2ff25a6c26e569a10.de-de.xlf"
$wsDeDe.Range("H6").Value = "2016-09-05 02:46:36"
$wsDeDe.Range("I6").Value = ""
$wsDeDe.Range("J6").Value = ""
$wsDeDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L6").Value = ""
$wsDeDe.Range("M6").Value = "True"
$wsDeDe.Range("N6").Value = ""
$wsDeDe.Range("O6").Value = "False"
$wsDeDe.Range("P6").Value = ""

$loDeDe.ListRows.Add()
$wsDeDe.Range("A7").Value = "953e0f08-aecb-4b53-a766-72bc9d9918db.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c21cabe021da961c62ae859f803b17101d2a2c0e/e2e/953e0f08-aecb-4b53-a766-72bc9d9918db.md", "", "", "953e0f08-aecb-4b53-a766-72bc9d9918db.md")
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e2e"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("F7").Value = "False"
$wsDeDe.Range("G7").Value = "953e0f08-aecb-4b53-a766-72bc9d9918db.c21cabe021da961c62ae859f803b17101d2a2c0e.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-09-05 02:46:36"
$wsDeDe.Range("I7").Value = ""
$wsDeDe.Range("J7").Value = ""
$wsDeDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L7").Value = ""
$wsDeDe.Range("M7").Value = "True"
$wsDeDe.Range("N7").Value = ""
$wsDeDe.Range("O7").Value = "False"
$wsDeDe.Range("P7").Value = ""
